$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new translation rows (64-81) introduced by this update
$ws.Cells.Item(64, 1).Value = 'ThingDef+VFES_Table_IlluminatedWardrobe.label'
$ws.Cells.Item(64, 2).Value = 'ThingDef'
$ws.Cells.Item(64, 3).Value = 'VFES_Table_IlluminatedWardrobe.label'
$ws.Cells.Item(64, 6).Value = '조명 대형 옷장'

$ws.Cells.Item(65, 1).Value = 'ThingDef+VFES_Table_IlluminatedWardrobe.description'
$ws.Cells.Item(65, 2).Value = 'ThingDef'
$ws.Cells.Item(65, 3).Value = 'VFES_Table_IlluminatedWardrobe.description'
$ws.Cells.Item(65, 6).Value = '부드럽게 빛나는 조명이 내장된 세련된 현대식 옷장입니다. 주변 침대의 편안함과 수면 효율을 약간 올려줍니다. 같은 침대 근처에 여러 개를 두어도 효과는 중첩되지 않습니다.'

$ws.Cells.Item(66, 1).Value = 'ThingDef+VFES_WallMountedTelevision.label'
$ws.Cells.Item(66, 2).Value = 'ThingDef'
$ws.Cells.Item(66, 3).Value = 'VFES_WallMountedTelevision.label'
$ws.Cells.Item(66, 6).Value = '벽걸이 TV'

$ws.Cells.Item(67, 1).Value = 'ThingDef+VFES_WallMountedTelevision.description'
$ws.Cells.Item(67, 2).Value = 'ThingDef'
$ws.Cells.Item(67, 3).Value = 'VFES_WallMountedTelevision.description'
$ws.Cells.Item(67, 6).Value = '얇은 고성능 평면 TV입니다. 공간을 적게 차지하면서도 일반 평면 TV와 동급의 선명한 화질과 색감을 제공합니다. 시야각은 넓지 않지만 벽에 설치할 수 있습니다.'

$ws.Cells.Item(68, 1).Value = 'ThingDef+VFES_TelevisionSpeaker.label'
$ws.Cells.Item(68, 2).Value = 'ThingDef'
$ws.Cells.Item(68, 3).Value = 'VFES_TelevisionSpeaker.label'
$ws.Cells.Item(68, 6).Value = 'TV 스피커'

$ws.Cells.Item(69, 1).Value = 'ThingDef+VFES_TelevisionSpeaker.description'
$ws.Cells.Item(69, 2).Value = 'ThingDef'
$ws.Cells.Item(69, 3).Value = 'VFES_TelevisionSpeaker.description'
$ws.Cells.Item(69, 6).Value = 'TV 시청 경험을 향상시키기 위해 설계된 고음질 오디오 시스템입니다. TV 근처에 배치하면 자동으로 연동되어 몰입감 있는 사운드를 제공하며, 시청 중인 정착민의 오락 욕구 충족도를 높입니다.'

$ws.Cells.Item(70, 1).Value = 'ThingDef+Spacer_OutdoorLamp.comps.3.offMessage'
$ws.Cells.Item(70, 2).Value = 'ThingDef'
$ws.Cells.Item(70, 3).Value = 'Spacer_OutdoorLamp.comps.3.offMessage'
$ws.Cells.Item(70, 6).Value = '전원 꺼짐: 충분한 햇빛 있음'

$ws.Cells.Item(71, 1).Value = 'ThingDef+VFES_AirPurifier.label'
$ws.Cells.Item(71, 2).Value = 'ThingDef'
$ws.Cells.Item(71, 3).Value = 'VFES_AirPurifier.label'
$ws.Cells.Item(71, 6).Value = '공기 청정기'

$ws.Cells.Item(72, 1).Value = 'ThingDef+VFES_AirPurifier.description'
$ws.Cells.Item(72, 2).Value = 'ThingDef'
$ws.Cells.Item(72, 3).Value = 'VFES_AirPurifier.description'
$ws.Cells.Item(72, 6).Value = '작고 효율적인 공기 청정기로, 먼지, 알레르기 유발 물질, 연기 등 유해 물질을 걸러내며 향기로운 공기를 순환시킵니다. 방의 청결도를 높이고 기분을 약간 향상시키며, 전반적인 환경을 개선합니다.'

$ws.Cells.Item(73, 1).Value = 'ThingDef+VFES_ModernArmchair.label'
$ws.Cells.Item(73, 2).Value = 'ThingDef'
$ws.Cells.Item(73, 3).Value = 'VFES_ModernArmchair.label'
$ws.Cells.Item(73, 6).Value = '현대식 팔걸이 의자'

$ws.Cells.Item(74, 1).Value = 'ThingDef+VFES_ModernArmchair.description'
$ws.Cells.Item(74, 2).Value = 'ThingDef'
$ws.Cells.Item(74, 3).Value = 'VFES_ModernArmchair.description'
$ws.Cells.Item(74, 6).Value = '우주 시대의 세련된 팔걸이 의자입니다. 메모리폼과 스마트 천으로 만들어졌으며 인체공학적 설계로 뛰어난 편안함과 미적 만족을 제공합니다.'

$ws.Cells.Item(75, 1).Value = 'ThingDef+VFES_ModernCouch.label'
$ws.Cells.Item(75, 2).Value = 'ThingDef'
$ws.Cells.Item(75, 3).Value = 'VFES_ModernCouch.label'
$ws.Cells.Item(75, 6).Value = '현대식 카우치'

$ws.Cells.Item(76, 1).Value = 'ThingDef+VFES_ModernCouch.description'
$ws.Cells.Item(76, 2).Value = 'ThingDef'
$ws.Cells.Item(76, 3).Value = 'VFES_ModernCouch.description'
$ws.Cells.Item(76, 6).Value = '우주 시대의 세련된 2인용 카우치입니다. 메모리폼과 인체공학적 설계, 스마트 천 덕에 최고의 편안함을 제공합니다.'

$ws.Cells.Item(77, 1).Value = 'ThoughtDef+VFES_FreshAir.stages.0.label'
$ws.Cells.Item(77, 2).Value = 'ThoughtDef'
$ws.Cells.Item(77, 3).Value = 'VFES_FreshAir.stages.0.label'
$ws.Cells.Item(77, 6).Value = '상쾌한 공기'

$ws.Cells.Item(78, 1).Value = 'ThoughtDef+VFES_FreshAir.stages.0.description'
$ws.Cells.Item(78, 2).Value = 'ThoughtDef'
$ws.Cells.Item(78, 3).Value = 'VFES_FreshAir.stages.0.description'
$ws.Cells.Item(78, 6).Value = '산뜻한 향이 나는 상쾌한 공기를 마셨어!'

$ws.Cells.Item(79, 1).Value = 'HediffDef+VFES_FreshAir.label'
$ws.Cells.Item(79, 2).Value = 'HediffDef'
$ws.Cells.Item(79, 3).Value = 'VFES_FreshAir.label'
$ws.Cells.Item(79, 6).Value = '상쾌한 공기'

$ws.Cells.Item(80, 1).Value = 'HediffDef+VFES_FreshAir.description'
$ws.Cells.Item(80, 2).Value = 'HediffDef'
$ws.Cells.Item(80, 3).Value = 'VFES_FreshAir.description'
$ws.Cells.Item(80, 6).Value = '산뜻한 향이 나는 상쾌한 공기를 마셨어!'

$ws.Cells.Item(81, 1).Value = 'JobDef+VFES_UseInteractiveTable.reportString'
$ws.Cells.Item(81, 2).Value = 'JobDef'
$ws.Cells.Item(81, 3).Value = 'VFES_UseInteractiveTable.reportString'
$ws.Cells.Item(81, 6).Value = '대화식 탁자 사용 중.'

# Restore the selection left by the author after editing the sheet
[void]$ws.Range("K83").Select()
